$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Professionalism")
$ws2 = $wb.Worksheets.Item("Hardware Development Process")

# Delete the old extra "Estimate Task Hours" row (old row 8) - rows below shift up by one
$ws2.Rows.Item(8).Delete()

# Copy the existing date-formatted cell style (m/d/yyyy) onto the Assigned/Due Date columns
# for the Architecture related rows, then fill in the real due dates.
$ws1.Range("C6").Copy() | Out-Null
$ws2.Range("C6:D10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws2.Range("C6").Value = (Get-Date -Year 2017 -Month 3 -Day 4).Date
$ws2.Range("D6").Value = (Get-Date -Year 2017 -Month 3 -Day 8).Date

$ws2.Range("C7").Value = (Get-Date -Year 2017 -Month 3 -Day 4).Date
$ws2.Range("D7").Value = (Get-Date -Year 2017 -Month 3 -Day 8).Date

$ws2.Range("C8").Value = (Get-Date -Year 2017 -Month 3 -Day 4).Date
$ws2.Range("D8").Value = (Get-Date -Year 2017 -Month 3 -Day 11).Date
$ws2.Range("E8").Value = "End of Class"

$ws2.Range("C9").Value = (Get-Date -Year 2017 -Month 3 -Day 4).Date
$ws2.Range("D9").Value = (Get-Date -Year 2017 -Month 3 -Day 12).Date
$ws2.Range("E9").Value = "End of Class"

$ws2.Range("C10").Value = (Get-Date -Year 2017 -Month 3 -Day 4).Date
$ws2.Range("D10").Value = (Get-Date -Year 2017 -Month 3 -Day 13).Date
$ws2.Range("E10").Value = "End of Class"

# Update view: selected cell moves to E11, no frozen/scrolled top-left cell
$ws2.Select()
$ws2.Range("E11").Select()

# Best-effort: update remembered window screen position (may not be persisted by this
# headless host, since there is no real on-screen window in this environment).
$wb.Windows.Item(1).Left = 3240
$wb.Windows.Item(1).Top = 460
